# Add a new worksheet "checkAboutHeaderAndVersionTest" after the existing
# "invalidCredentialTest" sheet, populate it with the About-page test data,
# and update the selections on both sheets.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new sheet right after the existing one so it becomes sheet2
# and the active tab.
$ws2 = $wb.Worksheets.Add([System.Type]::Missing, $ws1)
$ws2.Name = "checkAboutHeaderAndVersionTest"

# Header row
$ws2.Range("A1").Value = "username"
$ws2.Range("B1").Value = "password"
$ws2.Range("C1").Value = "language"
$ws2.Range("D1").Value = "expHeader"
$ws2.Range("E1").Value = "expVersion"

# Row 2 - admin
$ws2.Range("A2").Value = "admin"
$ws2.Range("B2").Value = "pass"
$ws2.Range("C2").Value = "English (Indian)"
$ws2.Range("D2").Value = "About OpenEMR"
$ws2.Range("E2").Value = "v6.0.0"

# Row 3 - accountant
$ws2.Range("A3").Value = "accountant"
$ws2.Range("B3").Value = "accountant"
$ws2.Range("C3").Value = "English (Indian)"
$ws2.Range("D3").Value = "About OpenEMR"
$ws2.Range("E3").Value = "v6.0.0"

# Restore the original sheet's selection (it loses tabSelected now that the
# new sheet is active) and set the new sheet's selection.
$ws1.Range("C2").Select()
$ws2.Range("B6").Select()
